$wb = $excel.ActiveWorkbook

function Fill-Row {
    param(
        $ws,
        [int]$r,
        [int]$aVal,
        [string]$b,
        [string]$c,
        [string]$d,
        [string]$e,
        $f,
        $g,
        [string]$h,
        [string]$i
    )
    $ws.Cells.Item($r, 1).Value = $aVal
    # Column B holds plain text dates (e.g. "2024-09-15"); force text format so it
    # is not silently re-interpreted as a date serial number.
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 9).Value = $i
}

function Update-Sheet {
    param($ws, [int]$lastRow)

    # Insert the new "合肥·ACGN夏日游园会第七回-泳池派对" row before the old row 3 (IE 动漫嘉年华).
    # Everything from row 3 onward shifts down by one row.
    $ws.Rows.Item(3).Insert()

    # Row 2: want-to-go count refreshed 689 -> 690
    $ws.Cells.Item(2, 6).Value = 690

    # New row 3: ACGN event
    Fill-Row $ws 3 2 "2024-09-15" "合肥·ACGN夏日游园会第七回-泳池派对" "金牛路金水里文化产业园 水善汇(金牛路店)" "2024.09.15 09:30-09.16 17:30" 2 50 "https://show.bilibili.com/platform/detail.html?id=91677" "//i2.hdslb.com/bfs/openplatform/202409/j9oW4hzR1725183897413.jpeg"

    # Renumber the sequential index column (A) for rows 4..(lastRow+1), which used to be rows 3..lastRow.
    for ($r = 4; $r -le ($lastRow + 1); $r++) {
        $ws.Cells.Item($r, 1).Value = ($r - 1)
    }

    # Refresh want-to-go counts that changed for events that only shifted rows.
    $ws.Cells.Item(5, 6).Value = 40     # 书香璃樱动漫游戏嘉年华: 39 -> 40
    $ws.Cells.Item(9, 6).Value = 3380   # 第十五届次元之门动漫游戏博览会: 3375 -> 3380
    $ws.Cells.Item(10, 6).Value = 4281  # 首届AT次元时代动漫游戏嘉年华: 4279 -> 4281

    # Insert the new "合肥·首届火影忍者同人only" row before the "W·A第五人格同人only2.0" row
    # (currently at row 11 after the first insert).
    $ws.Rows.Item(11).Insert()

    # New row 11: Naruto doujin event
    Fill-Row $ws 11 10 "2024-10-06" "合肥·首届火影忍者同人only" "长江东路金太阳家具广场南门二楼 优极篮球馆" "2024.10.06 09:30-10.06 17:30" 0 75 "https://show.bilibili.com/platform/detail.html?id=91658" "//i0.hdslb.com/bfs/openplatform/202408/f8ylbskH1725027552569.jpeg"

    # Renumber the sequential index column (A) for rows 12..(lastRow+2), which used to be rows 11..lastRow
    # (pre-first-insert numbering) / 10..(lastRow-1) originally.
    for ($r = 12; $r -le ($lastRow + 2); $r++) {
        $ws.Cells.Item($r, 1).Value = ($r - 1)
    }
}

# --- Sheet "展览" (exhibitions): originally rows 1..10 (A1:I10) ---
$ws1 = $wb.Worksheets.Item("展览")
Update-Sheet $ws1 10

# --- Sheet "全部类型" (all types): originally rows 1..13 (A1:I13) ---
$ws4 = $wb.Worksheets.Item("全部类型")
Update-Sheet $ws4 13
